# Update transition-probability matrix on Sheet1 with recalculated values
# reflecting games pulled March 7 (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1895043731778426
$ws.Range("C2").Value = 0.5685131195335277
$ws.Range("J2").Value = 0.008746355685131196
$ws.Range("P2").Value = 0.1516034985422741
$ws.Range("S2").Value = 0.08163265306122448
$ws.Range("B3").Value = 0.009708737864077669
$ws.Range("C3").Value = 0.04368932038834952
$ws.Range("J3").Value = 0.004854368932038835
$ws.Range("P3").Value = 0.8009708737864077
$ws.Range("S3").Value = 0.1407766990291262
$ws.Range("J4").Value = 0.04444444444444445
$ws.Range("P4").Value = 0.4444444444444444
$ws.Range("S4").Value = 0.5111111111111111
$ws.Range("B6").Value = 0.06511627906976744
$ws.Range("D6").Value = 0.004651162790697674
$ws.Range("F6").Value = 0.04186046511627907
$ws.Range("J6").Value = 0.2883720930232558
$ws.Range("O6").Value = 0.01395348837209302
$ws.Range("Q6").Value = 0.1162790697674419
$ws.Range("R6").Value = 0.1023255813953488
$ws.Range("S6").Value = 0.3674418604651163
$ws.Range("B7").Value = 0.1185567010309278
$ws.Range("D7").Value = 0.0154639175257732
$ws.Range("F7").Value = 0.04639175257731959
$ws.Range("J7").Value = 0.2061855670103093
$ws.Range("Q7").Value = 0.1701030927835052
$ws.Range("R7").Value = 0.07731958762886598
$ws.Range("S7").Value = 0.3659793814432989
$ws.Range("B8").Value = 0.1075794621026895
$ws.Range("D8").Value = 0.01466992665036675
$ws.Range("F8").Value = 0.0488997555012225
$ws.Range("J8").Value = 0.1124694376528117
$ws.Range("O8").Value = 0.009779951100244499
$ws.Range("Q8").Value = 0.1638141809290954
$ws.Range("R8").Value = 0.1271393643031785
$ws.Range("S8").Value = 0.4156479217603912
$ws.Range("B9").Value = 0.0847457627118644
$ws.Range("D9").Value = 0.01129943502824859
$ws.Range("F9").Value = 0.05649717514124294
$ws.Range("J9").Value = 0.0847457627118644
$ws.Range("O9").Value = 0.02259887005649718
$ws.Range("Q9").Value = 0.1751412429378531
$ws.Range("R9").Value = 0.1242937853107345
$ws.Range("S9").Value = 0.4406779661016949
$ws.Range("B10").Value = 0.1306607275426875
$ws.Range("D10").Value = 0.0244988864142539
$ws.Range("E10").Value = 0.0007423904974016332
$ws.Range("F10").Value = 0.07423904974016332
$ws.Range("J10").Value = 0.111358574610245
$ws.Range("O10").Value = 0.008166295471417966
$ws.Range("Q10").Value = 0.1907943578322197
$ws.Range("R10").Value = 0.09502598366740905
$ws.Range("S10").Value = 0.3645137342242019
$ws.Range("G11").Value = 0.145631067961165
$ws.Range("J11").Value = 0.0970873786407767
$ws.Range("K11").Value = 0.2200647249190938
$ws.Range("L11").Value = 0.4983818770226537
$ws.Range("S11").Value = 0.03883495145631068
$ws.Range("G12").Value = 0.73125
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.00625
$ws.Range("L12").Value = 0.01875
$ws.Range("S12").Value = 0.04375
$ws.Range("G13").Value = 0.6166666666666667
$ws.Range("J13").Value = 0.3166666666666667
$ws.Range("S13").Value = 0.06666666666666667
$ws.Range("F15").Value = 0.02898550724637681
$ws.Range("H15").Value = 0.1304347826086956
$ws.Range("I15").Value = 0.05797101449275362
$ws.Range("J15").Value = 0.4057971014492754
$ws.Range("K15").Value = 0.09178743961352658
$ws.Range("M15").Value = 0.01932367149758454
$ws.Range("O15").Value = 0.02415458937198068
$ws.Range("S15").Value = 0.2415458937198068
$ws.Range("F16").Value = 0.01809954751131222
$ws.Range("H16").Value = 0.1719457013574661
$ws.Range("I16").Value = 0.05429864253393665
$ws.Range("J16").Value = 0.4117647058823529
$ws.Range("K16").Value = 0.1266968325791855
$ws.Range("M16").Value = 0.02714932126696833
$ws.Range("O16").Value = 0.04977375565610859
$ws.Range("S16").Value = 0.1402714932126697
$ws.Range("F17").Value = 0.01686746987951807
$ws.Range("H17").Value = 0.1855421686746988
$ws.Range("I17").Value = 0.07710843373493977
$ws.Range("J17").Value = 0.4265060240963855
$ws.Range("K17").Value = 0.08674698795180723
$ws.Range("M17").Value = 0.02168674698795181
$ws.Range("N17").Value = 0.002409638554216868
$ws.Range("O17").Value = 0.05301204819277108
$ws.Range("S17").Value = 0.1301204819277108
$ws.Range("F18").Value = 0.02092050209205021
$ws.Range("H18").Value = 0.1297071129707113
$ws.Range("I18").Value = 0.07949790794979079
$ws.Range("J18").Value = 0.4435146443514644
$ws.Range("K18").Value = 0.06694560669456066
$ws.Range("M18").Value = 0.01255230125523013
$ws.Range("O18").Value = 0.1087866108786611
$ws.Range("S18").Value = 0.1380753138075314
$ws.Range("F19").Value = 0.0170940170940171
$ws.Range("H19").Value = 0.1864801864801865
$ws.Range("I19").Value = 0.07925407925407925
$ws.Range("J19").Value = 0.3869463869463869
$ws.Range("K19").Value = 0.1095571095571096
$ws.Range("M19").Value = 0.03108003108003108
$ws.Range("N19").Value = 0.005439005439005439
$ws.Range("O19").Value = 0.06837606837606838
$ws.Range("S19").Value = 0.1157731157731158
